# Emisión Motor - Varios vehiculos - General
# Arreglo en la emisión de varios vehiculos. Scripts de emisión de póliza
# complementaria y blanket.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Start from a clean slate: drop the old hyperlink + old data rows/styles
#    so we can rebuild the header + data block from scratch.
# ---------------------------------------------------------------------------
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("A1:Q3").ClearContents()
$ws.Range("A2:Q3").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Header row (row 1) - new column layout A:T
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value  = "idEjecucion"
$ws.Cells.Item(1,2).Value  = "Ambiente"
$ws.Cells.Item(1,3).Value  = "URL"
$ws.Cells.Item(1,4).Value  = "Usuario"
$ws.Cells.Item(1,5).Value  = "Contrasenia"
$ws.Cells.Item(1,6).Value  = "CodigoAgente"
$ws.Cells.Item(1,7).Value  = "NUM_GRUPO"
$ws.Cells.Item(1,8).Value  = "NroCuenta"
$ws.Cells.Item(1,9).Value  = "TIPOPOLIZA"
$ws.Cells.Item(1,10).Value = "CantVehiculos"
$ws.Cells.Item(1,11).Value = "TipoPlazo"
$ws.Cells.Item(1,12).Value = "MetodoDePago"
$ws.Cells.Item(1,13).Value = "OPCION_PAGOCUOTAS"
$ws.Cells.Item(1,14).Value = "CantCuotas"
$ws.Cells.Item(1,15).Value = "TipoTarjeta"
$ws.Cells.Item(1,16).Value = "NumTarjetaCred"
$ws.Cells.Item(1,17).Value = "FechaVencimiento"
$ws.Cells.Item(1,18).Value = "ConductoPago"
$ws.Cells.Item(1,19).Value = "NumCBU"
$ws.Cells.Item(1,20).Value = "FechaInicio"

# ---------------------------------------------------------------------------
# 3) Data rows 2-8
# ---------------------------------------------------------------------------
$idEjecucion   = @(2, 3, 4, 5, 6, 7, 8)
$codigoAgente  = @(1067, 6188, 1067, 4994, 234, 4994, 234)
$fechaInicio   = @("01/03/2020", "01/03/2020", "08/03/2020", "08/03/2020", "21/03/2020", "08/03/2020", "21/03/2020")

for ($i = 0; $i -lt 7; $i++) {
    $r = 2 + $i

    $ws.Cells.Item($r,1).Value  = $idEjecucion[$i]
    $ws.Cells.Item($r,2).Value  = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
    $ws.Cells.Item($r,3).Value  = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
    $ws.Cells.Item($r,4).Value  = "su"
    $ws.Cells.Item($r,5).Value  = "gw"
    $ws.Cells.Item($r,6).Value  = $codigoAgente[$i]
    $ws.Cells.Item($r,7).Value  = "Baioni Alejandro Luis"
    $ws.Cells.Item($r,8).Value  = 8892807402
    $ws.Cells.Item($r,9).Value  = "Motor"
    $ws.Cells.Item($r,10).Value = "Menos de 5 vehículos"
    $ws.Cells.Item($r,11).Value = "Anual"
    $ws.Cells.Item($r,12).Value = "Cupón"
    $ws.Cells.Item($r,13).Value = "No"
    $ws.Cells.Item($r,20).Value = "'" + $fechaInicio[$i]
}

# ---------------------------------------------------------------------------
# 4) Styles
#    - column B (Ambiente) and column T (FechaInicio) reuse the "quotePrefix"
#      plain style already present in the sheet (formerly used by A2/A3).
#    - column C (URL) reuses the hyperlink style (formerly used by B2/B3).
#    - column H (NroCuenta) reuses the numeric/grey-font style (formerly E2/E3).
#    - the date style for column T additionally carries a date NumberFormat.
# ---------------------------------------------------------------------------

# Give T2 its own quote-prefixed date style first (creates the single new
# cellXfs entry), then fan that exact style out to T3:T8 via copy/paste of
# formats only (so no further styles get created).
$ws.Cells.Item(2,20).NumberFormat = "mm-dd-yy"
$ws.Range("T2").Copy()
$ws.Range("T3:T8").PasteSpecial(-4122)

# Column B + T quote-prefix plain style (reuse existing style, originally on
# A2/A3 in the source workbook).
$ws.Range("A2").Copy()
$ws.Range("B2:B8").PasteSpecial(-4122)

# Column C hyperlink-like style (reuse existing style, originally on B2/B3).
$ws.Range("B2").Copy()
$ws.Range("C2:C8").PasteSpecial(-4122)

# Column H numeric/grey style (reuse existing style, originally on E2/E3).
$ws.Range("E2").Copy()
$ws.Range("H2:H8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Hyperlinks on column C (URL column)
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 7; $i++) {
    $r = 2 + $i
    $ws.Hyperlinks.Add($ws.Cells.Item($r,3), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do") | Out-Null
}

# ---------------------------------------------------------------------------
# 6) Column width for NUM_GRUPO (column G)
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 23.74

# ---------------------------------------------------------------------------
# 7) Page setup / view
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("G7").Select()
